$d = $word.ActiveDocument

# Remove the hidden "_GoBack" bookmark from the (currently empty) 4th
# paragraph, leaving the paragraph itself (and its formatting) intact.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# That paragraph is Paragraphs(4) -- insert the new TA-details paragraphs
# right after it, before the trailing blank paragraph.
$p = $d.Paragraphs(4)

$p1 = $p.Range.InsertParagraphAfter()
$d.Paragraphs(5).Range.Text = "Anil Kolla"

$p2 = $d.Paragraphs(5).Range.InsertParagraphAfter()
$d.Paragraphs(6).Range.Text = "Teaching Asssistant in OOP course at Northwest Missouri State University "

$p3 = $d.Paragraphs(6).Range.InsertParagraphAfter()
$d.Paragraphs(7).Range.Text = "United States of America"

$p4 = $d.Paragraphs(7).Range.InsertParagraphAfter()
$p5 = $d.Paragraphs(8).Range.InsertParagraphAfter()
